$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add I1 = "I0" and J1 = "IF", copying H1's formatting (style)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-8: I = 1, J = same value as H
for ($r = 2; $r -le 8; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}

$excel.CutCopyMode = 0
